# Map19Tto20T - fixing incorrect dimension item uids
#
# Every "Numerator" row in the 19Tto20TMap sheet had a distinct (and
# incorrect) dimension-item UID in column I ("num_or_den_uid"). They
# should all reference the single correct UID "Som9NRMQqV7" (the same
# one already used for the very first Numerator row). "Denominator"
# rows already reference the correct UID ("QpNj0nSuEhD") and are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("19Tto20TMap")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 8).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $h = $ws.Cells.Item($r, 8).Value2
    if ($h -eq "Numerator") {
        $ws.Cells.Item($r, 9).Value = "Som9NRMQqV7"
        # match the formatting already used by the correctly-mapped
        # Numerator rows (e.g. column G) rather than the old, buggy
        # font applied to the per-row unique uids
        $ws.Range("G" + $r).Copy()
        $ws.Cells.Item($r, 9).PasteSpecial(-4122) | Out-Null
    }
}
$excel.CutCopyMode = 0

# column I now only ever holds "Som9NRMQqV7" / "QpNj0nSuEhD" (both 11
# characters) instead of the old variable-length uids, so give it a
# tighter, fitted width
$ws.Columns.Item(9).ColumnWidth = 16.5

# reset the view: select I1, scrolled back to the left edge
$ws.Activate() | Out-Null
$ws.Range("I1").Select() | Out-Null
